$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: add hours worked (C20) and replace the placeholder "Deadline" description (D20)
# with the real release-notes entry. This also removes the now-unused "Deadline" shared
# string and appends the new description string (C24's SUM formula recalculates
# automatically from 30.5 to 39.5 once C20 is populated).
$ws.Cells.Item(20, 3).Value = 9
$ws.Cells.Item(20, 4).Value = "Major improvements in UI and logic. Twitter button and first release."

# Column D needs to grow to fit the new, longer description text.
$ws.Columns("D").ColumnWidth = 62.5

# Update the saved view: scroll so column B is left-most and select C21.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C21").Select() | Out-Null
